# Scheduled-runner update: refresh leve-profit inputs (currentAveragePrice /
# currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ) across the per-job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 450.6111
$ws.Range("I41").Value = 353.5
$ws.Range("K41").Value = 353.5
$ws.Range("M41").Value = 86.5
$ws.Range("H62").Value = 6000
$ws.Range("I62").Value = 3000
$ws.Range("K62").Value = 3000
$ws.Range("M62").Value = -2376
$ws.Range("H65").Value = 6000
$ws.Range("I65").Value = 3000
$ws.Range("K65").Value = 15000
$ws.Range("M65").Value = -11880
$ws.Range("H76").Value = 4192.35
$ws.Range("I76").Value = 3814.9375
$ws.Range("K76").Value = 3814.9375
$ws.Range("M76").Value = -3499.9375
$ws.Range("H79").Value = 4192.35
$ws.Range("I79").Value = 3814.9375
$ws.Range("K79").Value = 3814.9375
$ws.Range("M79").Value = -2722.9375
$ws.Range("H87").Value = 65500
$ws.Range("J87").Value = 65500
$ws.Range("L87").Value = 65500
$ws.Range("N87").Value = -67996
$ws.Range("H90").Value = 65500
$ws.Range("J90").Value = 65500
$ws.Range("L90").Value = 196500
$ws.Range("N90").Value = -208980
$ws.Range("H113").Value = 4166.25
$ws.Range("J113").Value = 4799.8
$ws.Range("L113").Value = 4799.8
$ws.Range("N113").Value = -11307.8
$ws.Range("H132").Value = 11465.474
$ws.Range("I132").Value = 12598.706
$ws.Range("K132").Value = 37796.118
$ws.Range("M132").Value = -35266.118
$ws.Range("H137").Value = 7142.5
$ws.Range("I137").Value = 2279.8823
$ws.Range("K137").Value = 6839.646900000001
$ws.Range("M137").Value = -4289.646900000001
$ws.Range("H138").Value = 2579.8
$ws.Range("I138").Value = 2495.2273
$ws.Range("K138").Value = 7485.6819
$ws.Range("M138").Value = -2345.6819
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2158.6155
$ws.Range("I45").Value = 1605.6364
$ws.Range("J45").Value = 5200
$ws.Range("K45").Value = 1605.6364
$ws.Range("L45").Value = 5200
$ws.Range("M45").Value = -1228.6364
$ws.Range("N45").Value = -5954
$ws.Range("H132").Value = 5338.8477
$ws.Range("I132").Value = 3459
$ws.Range("K132").Value = 10377
$ws.Range("M132").Value = -7847
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 10050
$ws.Range("J44").Value = 10050
$ws.Range("L44").Value = 10050
$ws.Range("N44").Value = -11044
$ws.Range("H86").Value = 2563.75
$ws.Range("I86").Value = 3075
$ws.Range("K86").Value = 3075
$ws.Range("M86").Value = -1952
$ws.Range("H89").Value = 2563.75
$ws.Range("I89").Value = 3075
$ws.Range("K89").Value = 15375
$ws.Range("M89").Value = -9759
$ws.Range("H105").Value = 5402.933
$ws.Range("I105").Value = 3294.0588
$ws.Range("J105").Value = 8160.6924
$ws.Range("K105").Value = 3294.0588
$ws.Range("L105").Value = 8160.6924
$ws.Range("M105").Value = -1547.0588
$ws.Range("N105").Value = -11654.6924
$ws.Range("H134").Value = 4360.2173
$ws.Range("I134").Value = 4734.3
$ws.Range("K134").Value = 14202.9
$ws.Range("M134").Value = -11667.9
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2028.6875
$ws.Range("I31").Value = 1875.7858
$ws.Range("K31").Value = 1875.7858
$ws.Range("M31").Value = -1580.7858
$ws.Range("H34").Value = 2028.6875
$ws.Range("I34").Value = 1875.7858
$ws.Range("K34").Value = 1875.7858
$ws.Range("M34").Value = -1673.7858
$ws.Range("H62").Value = 4998.6665
$ws.Range("I62").Value = 4998
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 4998
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -4374
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 4998.6665
$ws.Range("I65").Value = 4998
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 24990
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -21870
$ws.Range("N65").Value = -31240
$ws.Range("H94").Value = 7666.067
$ws.Range("I94").Value = 17032.834
$ws.Range("K94").Value = 17032.834
$ws.Range("M94").Value = -16581.834
$ws.Range("H107").Value = 852.9048
$ws.Range("I107").Value = 665.1875
$ws.Range("K107").Value = 665.1875
$ws.Range("M107").Value = 1254.8125
$ws.Range("H132").Value = 8376.467000000001
$ws.Range("I132").Value = 8520.482
$ws.Range("J132").Value = 4200
$ws.Range("K132").Value = 25561.446
$ws.Range("L132").Value = 12600
$ws.Range("M132").Value = -23031.446
$ws.Range("N132").Value = -17660
$ws.Range("H134").Value = 3695.9167
$ws.Range("I134").Value = 3625.1
$ws.Range("J134").Value = 4050
$ws.Range("K134").Value = 10875.3
$ws.Range("L134").Value = 12150
$ws.Range("M134").Value = -8340.299999999999
$ws.Range("N134").Value = -17220
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1782.1578
$ws.Range("I140").Value = 1603.3889
$ws.Range("K140").Value = 4810.1667
$ws.Range("M140").Value = 369.8333000000002
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 94.09524
$ws.Range("I2").Value = 51.666668
$ws.Range("J2").Value = 125.916664
$ws.Range("K2").Value = 51.666668
$ws.Range("L2").Value = 125.916664
$ws.Range("M2").Value = 61.333332
$ws.Range("N2").Value = -351.916664
$ws.Range("H113").Value = 1477
$ws.Range("I113").Value = 1474.7
$ws.Range("K113").Value = 1474.7
$ws.Range("M113").Value = 695.3
$ws.Range("H132").Value = 10664.115
$ws.Range("I132").Value = 9185.8125
$ws.Range("J132").Value = 13029.4
$ws.Range("K132").Value = 27557.4375
$ws.Range("L132").Value = 39088.2
$ws.Range("M132").Value = -25027.4375
$ws.Range("N132").Value = -44148.2
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 10161.728
$ws.Range("J40").Value = 6927.857
$ws.Range("L40").Value = 6927.857
$ws.Range("N40").Value = -7199.857
$ws.Range("H46").Value = 2824.6924
$ws.Range("I46").Value = 1281.8182
$ws.Range("K46").Value = 1281.8182
$ws.Range("M46").Value = -1093.8182
$ws.Range("H61").Value = 6536.0356
$ws.Range("I61").Value = 7188.6
$ws.Range("K61").Value = 7188.6
$ws.Range("M61").Value = -6986.6
$ws.Range("H99").Value = 44958
$ws.Range("I99").Value = 44947
$ws.Range("J99").Value = 44980
$ws.Range("K99").Value = 44947
$ws.Range("L99").Value = 44980
$ws.Range("M99").Value = -41952
$ws.Range("N99").Value = -50970
$ws.Range("H113").Value = 6536.0356
$ws.Range("I113").Value = 7188.6
$ws.Range("K113").Value = 7188.6
$ws.Range("M113").Value = -5018.6
$ws.Range("H132").Value = 4342.909
$ws.Range("I132").Value = 4342.909
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 13028.727
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -10498.727
$ws.Range("N132").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1620.675
$ws.Range("I132").Value = 1534.0513
$ws.Range("K132").Value = 4602.1539
$ws.Range("M132").Value = -2072.1539
